$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K (2022) data, mirroring the formatting pattern of column J
$ws.Cells.Item(4, 11).Value = 2022
$ws.Cells.Item(5, 11).Value = 3.9462868231169921
$ws.Cells.Item(6, 11).Value = 3.8007658934388928

# Copy the styles from column J so that K inherits compatible formatting,
# then apply the specific tweaks used by the new cellXfs entries.
$ws.Range("J4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("J5").Copy() | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null
$ws.Range("J6").Copy() | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the selection to match the target workbook state
$ws.Range("L5").Select() | Out-Null
